$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '45.513.39'
$ws.Range("E2").Value = '  +7.28%  '
$ws.Range("D3").Value = '2.378.61'
$ws.Range("E3").Value = '  +4.27%  '
$ws.Range("E4").Value = '  +0.45%  '
$ws.Range("B5").Value = 'Solana'
$ws.Range("C5").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '111.48'
$ws.Range("E5").Value = '  +7.55%  '
$ws.Range("B6").Value = 'BNB'
$ws.Range("C6").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '317.19'
$ws.Range("E6").Value = '  +2.06%  '
$ws.Range("E7").Value = '  +3.80%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.630'
$ws.Range("E9").Value = '  +5.28%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.93'
$ws.Range("E10").Value = '  +7.97%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0929'
$ws.Range("E11").Value = '  +3.39%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.65'
$ws.Range("E12").Value = '  +5.59%  '
$ws.Range("E13").Value = '  +4.79%  '
$ws.Range("E14").Value = '  +0.86%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.74'
$ws.Range("E15").Value = '  +4.76%  '
$ws.Range("D16").Value = '2.739.25'
$ws.Range("E16").Value = '  +4.17%  '
$ws.Range("D17").Value = '2.394.01'
$ws.Range("E17").Value = '  +5.34%  '
$ws.Range("D18").Value = '45.289.16'
$ws.Range("E18").Value = '  +6.93%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.63'
$ws.Range("E19").Value = '  +5.42%  '
$ws.Range("E20").Value = '  +4.25%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.04'
$ws.Range("E21").Value = '  -3.11%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '75.07'
$ws.Range("E22").Value = '  +3.12%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.54'
$ws.Range("E23").Value = '  +3.91%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '268.31'
$ws.Range("E24").Value = '  +2.00%  '
$ws.Range("E25").Value = '  +7.57%  '
$ws.Range("E26").Value = '  -0.59%  '
$ws.Range("E27").Value = '  +8.55%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.32'
$ws.Range("E28").Value = '  +6.27%  '
$ws.Range("E29").Value = '  +3.03%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.89'
$ws.Range("E30").Value = '  +3.04%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '38.68'
$ws.Range("E31").Value = '  +8.40%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0941'
$ws.Range("E32").Value = '  +9.97%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '169.95'
$ws.Range("E33").Value = '  +3.34%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.04'
$ws.Range("E34").Value = '  +18.62%  '
$ws.Range("E35").Value = '  +2.68%  '
$ws.Range("B36").Value = 'RenderToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.87'
$ws.Range("E36").Value = '  +8.71%  '
$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.117'
$ws.Range("E37").Value = '  +4.76%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.06'
$ws.Range("E38").Value = '  +13.11%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0366'
$ws.Range("E39").Value = '  +5.13%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.92'
$ws.Range("E40").Value = '  +5.55%  '
$ws.Range("E41").Value = '  +11.40%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '106.41'
$ws.Range("E42").Value = '  +8.25%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '13.83'
$ws.Range("E43").Value = '  +16.36%  '
$ws.Range("E44").Value = '  +6.52%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '71.73'
$ws.Range("E45").Value = '  +4.22%  '
$ws.Range("E46").Value = '  +0.20%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '118.02'
$ws.Range("E47").Value = '  +7.81%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.78'
$ws.Range("E48").Value = '  +12.07%  '
$ws.Range("E49").Value = '  +19.77%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '79.28'
$ws.Range("E50").Value = '  +2.89%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.18'
$ws.Range("E51").Value = '  +6.69%  '
